$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1 - Product Burndown")

$ws.Range("E14").Value = 15
$ws.Range("B15").Value = 12

$ws.Range("D16").Select()
